$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15119
$ws1.Range("F3").Value = 19365
$ws1.Range("F4").Value = 148
$ws1.Range("F5").Value = 156
$ws1.Range("F14").Value = 193
$ws1.Range("F15").Value = 238
$ws1.Range("F20").Value = 105
$ws1.Range("F21").Value = 243
$ws1.Range("F22").Value = 8113
$ws1.Range("F24").Value = 40
$ws1.Range("F25").Value = 8
$ws1.Range("F27").Value = 1265
$ws1.Range("F29").Value = 11
$ws1.Range("F31").Value = 6505
$ws1.Range("F32").Value = 127
$ws1.Range("F34").Value = 181
$ws1.Range("F36").Value = 297
$ws1.Range("F37").Value = 5523
$ws1.Range("F38").Value = 1013
$ws1.Range("F41").Value = 58

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15119
$ws4.Range("F3").Value = 19365
$ws4.Range("F4").Value = 148
$ws4.Range("F5").Value = 156
$ws4.Range("F14").Value = 193
$ws4.Range("F15").Value = 238
$ws4.Range("F21").Value = 105
$ws4.Range("F22").Value = 243
$ws4.Range("F23").Value = 8113
$ws4.Range("F25").Value = 40
$ws4.Range("F26").Value = 8
$ws4.Range("F28").Value = 1265
$ws4.Range("F30").Value = 11
$ws4.Range("F34").Value = 6505
$ws4.Range("F35").Value = 127
$ws4.Range("F37").Value = 181
$ws4.Range("F39").Value = 297
$ws4.Range("F40").Value = 5523
$ws4.Range("F41").Value = 1013
$ws4.Range("F44").Value = 58
